# --- edit.ps1 ---
# Applies the "update from richard nov 23" revision to wages_historical.xlsx:
#   1. Break the stale external link to wages_historical.xlsx (FRED/CPI pull)
#      so every `=[1]FRED!A..` formula in column A collapses to its cached
#      literal value and the externalLinks parts + workbook.xml reference
#      are dropped on save.
#   2. Bump the report title from "... September 2015" to "... October 2015".
#   3. Refresh the Y/Y wage-growth figures for the existing last data row (98).
#   4. Append the new trailing data row (99) for the new month.
#   5. Move the live selection to C92, as in the source workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Break the external workbook link (wages_historical.xlsx).
$wb.BreakLink("wages_historical.xlsx", 1)

# 2. Update the report title (row 1) - "September 2015" -> "October 2015".
$ws.Range("A1").Value2 = "Real Average Weekly Wages; Year/Year: January 2008 -- October 2015"

# 3. Row 98 (period ending 2015-09-01): refreshed Y/Y wage-growth figures for
#    columns C:BA (A and B are unchanged).
$ws.Range("C98").Value2 = [double]"1.1005017013327831"
$ws.Range("D98").Value2 = [double]"2.1306278272832739"
$ws.Range("E98").Value2 = [double]"3.4696053285147563"
$ws.Range("F98").Value2 = [double]"-1.4526918822004753"
$ws.Range("G98").Value2 = [double]"1.894983351928399"
$ws.Range("H98").Value2 = [double]"0.74700110289772503"
$ws.Range("I98").Value2 = [double]"2.8513267917191158"
$ws.Range("J98").Value2 = [double]"6.1646266236050664"
$ws.Range("K98").Value2 = [double]"-5.684904347437711"
$ws.Range("L98").Value2 = [double]"1.4497611988441521"
$ws.Range("M98").Value2 = [double]"-3.1705414087965629E-2"
$ws.Range("N98").Value2 = [double]"-3.4426260724629052"
$ws.Range("O98").Value2 = [double]"1.673531816142465"
$ws.Range("P98").Value2 = [double]"1.6667312468045905"
$ws.Range("Q98").Value2 = [double]"0.71629687102618744"
$ws.Range("R98").Value2 = [double]"3.1267199055027248"
$ws.Range("S98").Value2 = [double]"7.1290253060657491E-2"
$ws.Range("T98").Value2 = [double]"4.563695160672947"
$ws.Range("U98").Value2 = [double]"-1.9086364013446206"
$ws.Range("V98").Value2 = [double]"2.7999722138947103"
$ws.Range("W98").Value2 = [double]"-4.334800042173066E-3"
$ws.Range("X98").Value2 = [double]"2.9855764779201559"
$ws.Range("Y98").Value2 = [double]"1.8543778500428487"
$ws.Range("Z98").Value2 = [double]"0.41437712885608741"
$ws.Range("AA98").Value2 = [double]"-2.6519795651680766"
$ws.Range("AB98").Value2 = [double]"8.1967882980174836E-2"
$ws.Range("AC98").Value2 = [double]"2.5326894631134462"
$ws.Range("AD98").Value2 = [double]"5.0756568510671753"
$ws.Range("AE98").Value2 = [double]"6.4064024621575291"
$ws.Range("AF98").Value2 = [double]"2.6559036258003519"
$ws.Range("AG98").Value2 = [double]"3.2279932806311606"
$ws.Range("AH98").Value2 = [double]"-0.59087405653857894"
$ws.Range("AI98").Value2 = [double]"1.5170276511489749"
$ws.Range("AJ98").Value2 = [double]"1.1924621450056274"
$ws.Range("AK98").Value2 = [double]"-0.17736380203700494"
$ws.Range("AL98").Value2 = [double]"2.5440748920622096"
$ws.Range("AM98").Value2 = [double]"-1.2932586777004671"
$ws.Range("AN98").Value2 = [double]"2.7497568290862908"
$ws.Range("AO98").Value2 = [double]"2.9095014233026415"
$ws.Range("AP98").Value2 = [double]"0.70025673522697107"
$ws.Range("AQ98").Value2 = [double]"1.7281385442674255"
$ws.Range("AR98").Value2 = [double]"3.5340516634853874"
$ws.Range("AS98").Value2 = [double]"0.90683564915752368"
$ws.Range("AT98").Value2 = [double]"-0.38870022789273728"
$ws.Range("AU98").Value2 = [double]"1.0454474201444657"
$ws.Range("AV98").Value2 = [double]"2.3737022812335042"
$ws.Range("AW98").Value2 = [double]"4.6355056937479464"
$ws.Range("AX98").Value2 = [double]"4.30470046454833"
$ws.Range("AY98").Value2 = [double]"-0.66096357107067227"
$ws.Range("AZ98").Value2 = [double]"-0.27867177946022642"
$ws.Range("BA98").Value2 = [double]"-3.4194571934963092"

# 4. Row 99 (period ending 2015-10-01): brand new data row appended to the
#    table. Match number formatting of the row above (style s="6" / s="5")
#    before writing values, since a never-before-used cell otherwise falls
#    back to the bare column style.
$ws.Range("A99:BA99").NumberFormat = $ws.Range("A98:BA98").NumberFormat
$ws.Range("A99").Value2 = [double]"42278"
$ws.Range("B99").Value2 = [double]"2.0109720293641753"
$ws.Range("C99").Value2 = [double]"3.875389073312093"
$ws.Range("D99").Value2 = [double]"1.6713264370977241"
$ws.Range("E99").Value2 = [double]"3.4988375535301701"
$ws.Range("F99").Value2 = [double]"0.43302394287164692"
$ws.Range("G99").Value2 = [double]"2.0792631730655082"
$ws.Range("H99").Value2 = [double]"1.2589464957269563"
$ws.Range("I99").Value2 = [double]"3.5758594401086343"
$ws.Range("J99").Value2 = [double]"6.3838400348970579"
$ws.Range("K99").Value2 = [double]"-8.3058752488200565"
$ws.Range("L99").Value2 = [double]"2.6486635728379273"
$ws.Range("M99").Value2 = [double]"1.6996638775173722"
$ws.Range("N99").Value2 = [double]"-1.4022928484078947"
$ws.Range("O99").Value2 = [double]"1.7020244649269953"
$ws.Range("P99").Value2 = [double]"1.6524800621482143"
$ws.Range("Q99").Value2 = [double]"1.0015035399525174"
$ws.Range("R99").Value2 = [double]"4.1900086022102334"
$ws.Range("S99").Value2 = [double]"0.70079406355546436"
$ws.Range("T99").Value2 = [double]"4.4161672817386926"
$ws.Range("U99").Value2 = [double]"-0.65807158333573323"
$ws.Range("V99").Value2 = [double]"3.1233671692587106"
$ws.Range("W99").Value2 = [double]"1.559358402569391"
$ws.Range("X99").Value2 = [double]"3.7401779270501736"
$ws.Range("Y99").Value2 = [double]"2.7034601909267413"
$ws.Range("Z99").Value2 = [double]"2.4437636772824591"
$ws.Range("AA99").Value2 = [double]"-0.95244836323532156"
$ws.Range("AB99").Value2 = [double]"-0.44515591527256249"
$ws.Range("AC99").Value2 = [double]"1.485635963187746"
$ws.Range("AD99").Value2 = [double]"6.0250952034602738"
$ws.Range("AE99").Value2 = [double]"5.8004642924436753"
$ws.Range("AF99").Value2 = [double]"4.202902598420188"
$ws.Range("AG99").Value2 = [double]"3.161060117917744"
$ws.Range("AH99").Value2 = [double]"-1.8565824674211004"
$ws.Range("AI99").Value2 = [double]"2.077786971076335"
$ws.Range("AJ99").Value2 = [double]"3.1326224220425574"
$ws.Range("AK99").Value2 = [double]"0.37679809409480747"
$ws.Range("AL99").Value2 = [double]"2.5351355087312619"
$ws.Range("AM99").Value2 = [double]"-0.71203965732984431"
$ws.Range("AN99").Value2 = [double]"3.3543683199632057"
$ws.Range("AO99").Value2 = [double]"3.6929845617981316"
$ws.Range("AP99").Value2 = [double]"-0.2734690259049502"
$ws.Range("AQ99").Value2 = [double]"2.2594997596382607"
$ws.Range("AR99").Value2 = [double]"4.6245140590792895"
$ws.Range("AS99").Value2 = [double]"2.1989038850612483"
$ws.Range("AT99").Value2 = [double]"-0.31218107179246302"
$ws.Range("AU99").Value2 = [double]"1.044125285500326"
$ws.Range("AV99").Value2 = [double]"3.0328898921944356"
$ws.Range("AW99").Value2 = [double]"3.4764430084489697"
$ws.Range("AX99").Value2 = [double]"5.6330882589008358"
$ws.Range("AY99").Value2 = [double]"0.4140175531080646"
$ws.Range("AZ99").Value2 = [double]"0.99756798909318467"
$ws.Range("BA99").Value2 = [double]"-2.8344035384745156"

# 5. Window state: move the selection to C92 (matches the diff's <selection>).
$ws.Range("C92").Select()

